# Populate the new "Name.Number" column (E) with fixed (previously randomly
# generated) values, per commit "Removed random number generation for
# nameNumber to keep test data consistent".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("E1").Value = 'Name.Number'

# Data rows E2:E225 - last name + "." + number, some values carry a
# leading space (preserved verbatim from the source data).
$nameNumbers = @(
    'Koop.960',
    ' Mcbain.358',
    ' Huot.633',
    ' Bezanson.660',
    ' Cavalieri.839',
    ' Bass.111',
    ' Isaacs.442',
    ' Wolanski.136',
    ' Clendening.100',
    ' Ashbaugh.213',
    ' Pinder.964',
    ' Depalma.763',
    ' Imel.715',
    ' Shambaugh.40',
    ' Albritton.668',
    ' Slocum.139',
    ' Maysonet.58',
    ' Maag.400',
    ' Laguardia.605',
    ' Rosso.903',
    ' Jeanbaptiste.831',
    ' Siegrist.361',
    ' Balas.71',
    ' Costello.760',
    ' Poulter.816',
    ' Whitmire.337',
    ' Politte.833',
    ' Hartshorn.948',
    ' Iser.962',
    ' Tinajero.70',
    ' Caplan.540',
    ' Seppala.173',
    ' Statler.8',
    ' Dartez.850',
    ' Backman.359',
    ' Norgard.983',
    ' Ingwersen.694',
    'Torgersen.942',
    ' Landey.624',
    ' Herbst.123',
    ' Mings.349',
    ' Gandara.979',
    ' Caffrey.986',
    ' Branning.888',
    ' Beland.117',
    ' Mariscal.729',
    ' Medel.318',
    ' Normandin.386',
    ' Breuer.491',
    ' Liberto.857',
    ' Godfrey.949',
    ' Mikkelson.973',
    ' Fleagle.381',
    ' Darlington.481',
    ' Hosein.882',
    ' Juarbe.11',
    ' Deslauriers.767',
    ' Horney.105',
    ' Nivens.576',
    ' Enfinger.107',
    ' Mccraw.316',
    ' Eggleton.598',
    ' Wymore.153',
    ' Dayton.843',
    ' Petti.427',
    ' Guercio.396',
    ' Vasconcellos.50',
    ' Heal.856',
    ' Hornick.940',
    ' Franson.807',
    ' Mahar.143',
    ' Denney.753',
    ' Saffell.584',
    ' Maze.620',
    ' Pawlowicz.398',
    ' Lahr.942',
    ' Lerman.315',
    'Provencal.801',
    ' Delagarza.200',
    ' Thompkins.432',
    ' Mccroskey.8',
    ' Wynter.419',
    ' Cliff.333',
    ' Storck.591',
    ' Meikle.556',
    ' Lightfoot.559',
    ' Ester.921',
    ' Mungo.80',
    ' Michalski.556',
    ' Lunceford.641',
    ' Siebert.643',
    ' Call.939',
    ' Robin.11',
    ' Roddy.883',
    ' Pangle.362',
    ' Hastings.286',
    ' Seller.821',
    ' Mumaw.240',
    ' Popham.241',
    ' Grieve.693',
    ' Camacho.4',
    ' Ehrenberg.441',
    ' Champine.202',
    ' Wells.55',
    ' Spicer.41',
    ' Shirley.671',
    ' Eddington.306',
    ' Morrisette.605',
    ' Whiting.851',
    ' Nordahl.9',
    ' Mullet.321',
    ' Rand.840',
    ' Mathisen.218',
    ' Garten.560',
    ' Tyree.611',
    ' Peppler.143',
    ' Paulson.371',
    'Mercer.633',
    ' Lickteig.332',
    ' Haddon.17',
    ' Mauger.736',
    ' Berge.171',
    ' Willcutt.973',
    ' Gaver.489',
    ' Ruyle.883',
    ' Heeter.44',
    ' Fitzwater.145',
    ' Aquilino.382',
    ' Palomino.241',
    ' Dunson.706',
    ' Hager.952',
    ' Bussell.732',
    ' Reese.274',
    ' Borjas.930',
    ' Cady.595',
    ' Okada.299',
    ' Debow.528',
    ' Calderone.593',
    ' Berrier.471',
    ' Parry.560',
    ' Tignor.161',
    ' Cutsforth.174',
    ' Cripe.236',
    ' Belmont.982',
    ' Girard.502',
    ' Vining.935',
    ' Clever.468',
    ' Dias.832',
    ' Drake.377',
    ' Lane.218',
    ' Gerry.184',
    ' Younts.552',
    ' Mcmurtrie.820',
    ' Dowless.286',
    ' Scheffer.323',
    ' Rempel.454',
    ' Lis.287',
    'Marsch.502',
    ' Lesperance.914',
    ' Gunning.855',
    ' Cripps.973',
    ' Younkin.677',
    ' Braz.539',
    ' Zielinski.230',
    ' Junior.195',
    ' Veneziano.661',
    ' Down.797',
    ' Ogle.891',
    ' Resnick.437',
    ' Tindell.641',
    ' Woltz.563',
    ' Borjas.78',
    ' Esquer.919',
    ' Zdenek.728',
    ' Astorga.32',
    ' Dollins.210',
    ' Standifer.77',
    ' Stoke.982',
    ' Grief.860',
    ' Denis.806',
    ' Gamblin.788',
    ' Faddis.634',
    ' Binns.159',
    ' Hepfer.154',
    ' Slade.484',
    ' Papp.563',
    ' Encinas.262',
    ' Mcdonald.342',
    ' Benigno.116',
    ' Brinkley.866',
    ' Certain.335',
    ' Cypher.353',
    ' Lautenschlage.917',
    ' Rushford.890',
    ' Delossantos.38',
    ' Slee.951',
    ' Moretti.365',
    'Alphin.603',
    ' Rainer.108',
    ' Breitenstein.677',
    ' Coletta.344',
    ' Hieber.798',
    ' Closson.956',
    ' Tully.114',
    ' Caywood.170',
    ' Silsby.242',
    ' Tait.159',
    ' Aleman.73',
    ' Brindley.439',
    ' Dicarlo.674',
    ' Florez.613',
    ' Sawyers.890',
    ' Byrne.461',
    ' Hagar.119',
    ' Goldie.250',
    ' Monfort.301',
    ' Toll.679',
    ' Gammage.15',
    ' Darnell.242',
    ' Paro.334',
    ' Barwick.128',
    ' Dugger.786',
    ' Leibowitz.835',
    ' Krebsbach.938'
)

for ($i = 0; $i -lt $nameNumbers.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 5).Value = $nameNumbers[$i]
}

# Re-select the full column E (matches the saved selection state).
$ws.Range("E1:E1048576").Select()

# The saved workbook also carries one extra (unused) conditional-format
# "duplicate values" dxf entry in styles.xml, identical to the pair that
# was already there protecting E2. Recreate then remove a matching rule
# on a scratch cell so its dxf is appended without leaving behind any
# new conditional-formatting rule or touching the used range/dimension.
$scratch = $ws.Range("Z1")
$fc = $scratch.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = -16383844
$fc.Interior.Color = 13551615
$scratch.FormatConditions.Delete()
